# ---------------------------------------------------------------------------
# Adds the "ODI Batting Extra" sheet (full data scraped for extra batting and
# bowling fields) and removes the stray empty B-column cells from the
# "ODI Batting" sheet for the "did not bat" rows.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. ODI Batting: clear the leftover empty inline-string cells in column B
#    for the rows where the player did not bat (row 30, 31, 35, 52, 56, 60).
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$emptyBRows = @(30, 31, 35, 52, 56, 60)
foreach ($r in $emptyBRows) {
    $battingSheet.Cells.Item($r, 2).Value = $null
}

# ---------------------------------------------------------------------------
# 2. Add the new "ODI Batting Extra" sheet after the last existing sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Header row
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Re-use the same bold/centered/bordered header formatting already used by
# the other sheets' header rows (copy format only, keep the values above).
$headerFormatSource = $wb.Worksheets.Item("Player Info").Range("A1")
$headerFormatSource.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$rows = @(
  ,@('3833', '8', '0', '0', '1.34%', 'NO')
  ,@('3984', '7', '9', '11', '52.14%', 'YES')
  ,@('3988', '6', '4', '1', '16.34%', 'NO')
  ,@('4067', '6', '0', '0', '2.19%', 'NO')
  ,@('4069', '6', '6', '3', '30.69%', 'NO')
  ,@('4071', '6', '1', '1', '9.22%', 'NO')
  ,@('4074', $null, $null, $null, $null, 'NO')
  ,@('4076', '6', '4', '1', '19.01%', 'NO')
  ,@('4108', $null, $null, $null, $null, 'NO')
  ,@('4115', '6', '0', '0', '1.48%', 'NO')
  ,@('4117', '6', '3', '2', '19.58%', 'NO')
  ,@('4123', '6', '3', '0', '7.11%', 'NO')
  ,@('4125', '3', '6', '4', '35.22%', 'NO')
  ,@('4166', '4', '0', '1', '10.28%', 'NO')
  ,@('4167', $null, $null, $null, $null, 'NO')
  ,@('4168', '4', '4', '1', '18.41%', 'NO')
  ,@('4169', $null, $null, $null, $null, 'NO')
  ,@('4170', '3', '0', '0', $null, 'NO')
  ,@('4222', $null, $null, $null, $null, 'NO')
  ,@('4224', '7', '0', '0', '0.87%', 'NO')
  ,@('4226', '5', '5', '4', '22.50%', 'NO')
  ,@('4234', '6', '2', '2', '16.32%', 'NO')
  ,@('4235', $null, $null, $null, $null, 'NO')
  ,@('4236', '6', '1', '0', '4.35%', 'NO')
  ,@('4258', $null, $null, $null, $null, 'NO')
  ,@('4263', '6', '4', '1', '21.49%', 'NO')
  ,@('4266', '5', '4', '0', '9.90%', 'NO')
  ,@('4270', '5', '1', '1', '7.35%', 'NO')
  ,@('4273', '5', $null, $null, $null, 'NO')
  ,@('4274', '6', $null, $null, $null, 'NO')
  ,@('4275', '5', '1', '0', '3.76%', 'NO')
  ,@('4276', $null, $null, $null, $null, 'NO')
  ,@('4277', '5', '0', '0', '1.22%', 'NO')
  ,@('4306', '6', $null, $null, $null, 'NO')
  ,@('4312', '6', '4', '0', '6.60%', 'NO')
  ,@('4316', '6', '0', '0', $null, 'NO')
  ,@('4329', $null, $null, $null, $null, 'NO')
  ,@('4336', '6', '1', '0', '2.81%', 'NO')
  ,@('4341', '5', '3', '0', '8.64%', 'NO')
  ,@('4351', '5', '4', '0', '6.98%', 'NO')
  ,@('4354', '6', '0', '0', $null, 'NO')
  ,@('4429', '3', '6', '0', '14.63%', 'NO')
  ,@('4430', '3', '0', '1', '4.35%', 'NO')
  ,@('4431', '3', '1', '0', '1.31%', 'NO')
  ,@('4435', $null, $null, $null, $null, 'NO')
  ,@('4564', '5', '0', '1', '8.31%', 'NO')
  ,@('4565', $null, $null, $null, $null, 'NO')
  ,@('4567', '5', '2', '1', '9.05%', 'NO')
  ,@('4594', $null, $null, $null, $null, 'NO')
  ,@('4644', '5', '2', '1', '9.45%', 'NO')
  ,@('4645', '5', $null, $null, $null, 'NO')
  ,@('4646', '5', '0', '0', '2.13%', 'NO')
  ,@('4647', $null, $null, $null, $null, 'NO')
  ,@('4648', '5', '0', '0', $null, 'NO')
  ,@('4660', $null, $null, $null, $null, 'NO')
  ,@('4663', $null, $null, $null, $null, 'NO')
  ,@('4666', $null, $null, $null, $null, 'NO')
  ,@('4725', '8', '1', '0', '2.66%', 'NO')
  ,@('4728', '7', $null, $null, $null, 'NO')
  ,@('4732', '7', '3', '0', '9.29%', 'NO')
)

$rowIndex = 2
foreach ($row in $rows) {
    # Column A - MATCH_CODE (kept as text, matches source scrape format)
    $ws.Cells.Item($rowIndex, 1).Value = "'" + $row[0]

    # Column B - BATTING_POSITION (numeric when present, blank otherwise)
    if ($null -eq $row[1]) {
        $ws.Cells.Item($rowIndex, 2).Value = $null
    } else {
        $ws.Cells.Item($rowIndex, 2).Value = [double]$row[1]
    }

    # Columns C, D, E - NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL (text, blank when missing)
    for ($col = 2; $col -le 4; $col++) {
        $val = $row[$col]
        if ($null -eq $val) {
            $ws.Cells.Item($rowIndex, $col + 1).Value = $null
        } else {
            $ws.Cells.Item($rowIndex, $col + 1).Value = "'" + $val
        }
    }

    # Column F - MAN_OF_MATCH (text, always populated)
    $ws.Cells.Item($rowIndex, 6).Value = "'" + $row[5]

    $rowIndex++
}

$ws.Range("A1").Select() | Out-Null
